$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.774.52'
$ws.Range('E2').Value = '  -0.86%  '
$ws.Range('D3').Value = '1.598.07'
$ws.Range('E3').Value = '  -2.44%  '
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').Value = "'208.66"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -2.73%  '
$ws.Range('E6').Value = '  -0.05%  '
$ws.Range('E7').Value = '  -4.83%  '
$ws.Range('D8').Value = "'0.246"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -2.36%  '
$ws.Range('D9').Value = "'0.0609"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -2.45%  '
$ws.Range('D10').Value = "'17.85"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -3.60%  '
$ws.Range('E11').Value = '  -0.71%  '
$ws.Range('D12').Value = '1.822.29'
$ws.Range('E12').Value = '  -2.24%  '
$ws.Range('D13').Value = '1.584.84'
$ws.Range('E13').Value = '  -3.56%  '
$ws.Range('D14').Value = "'4.04"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -4.02%  '
$ws.Range('D15').Value = "'0.508"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -4.25%  '
$ws.Range('D16').Value = '25.775.57'
$ws.Range('E16').Value = '  -0.85%  '
$ws.Range('D17').Value = "'60.38"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -2.39%  '
$ws.Range('D18').Value = '0.0₃0714'
$ws.Range('E18').Value = '  -4.20%  '
$ws.Range('D20').Value = "'189.19"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.08%  '
$ws.Range('D21').Value = "'4.17"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.75%  '
$ws.Range('D22').Value = "'9.32"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -4.26%  '
$ws.Range('D23').Value = "'5.93"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -3.08%  '
$ws.Range('B24').Value = 'BinanceUSD'
$ws.Range('C24').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D24').Value = "'1.01"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.05%  '
$ws.Range('B25').Value = 'Stellar'
$ws.Range('C25').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D25').Value = "'0.128"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -3.88%  '
$ws.Range('D26').Value = "'140.84"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -2.09%  '
$ws.Range('E27').Value = '  -4.76%  '
$ws.Range('B28').Value = 'EthereumClassic'
$ws.Range('C28').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D28').Value = "'14.97"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -2.00%  '
$ws.Range('B29').Value = 'Cosmos'
$ws.Range('C29').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D29').Value = "'6.50"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -5.19%  '
$ws.Range('E30').Value = '  -3.95%  '
$ws.Range('E31').Value = '  -3.14%  '
$ws.Range('D32').Value = "'3.06"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -3.05%  '
$ws.Range('E33').Value = '  -4.88%  '
$ws.Range('D34').Value = "'2.40"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.46%  '
$ws.Range('E35').Value = '  -1.78%  '
$ws.Range('D36').Value = '1.093.70'
$ws.Range('E36').Value = '  -3.70%  '
$ws.Range('D37').Value = "'2.37"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -2.88%  '
$ws.Range('E38').Value = '  -0.47%  '
$ws.Range('E39').Value = '  -2.36%  '
$ws.Range('D40').Value = "'0.789"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -9.16%  '
$ws.Range('D41').Value = "'0.496"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -4.72%  '
$ws.Range('D42').Value = "'95.35"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -3.11%  '
$ws.Range('D43').Value = '1.735.27'
$ws.Range('E43').Value = '  -2.15%  '
$ws.Range('D44').Value = "'5.06"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -3.31%  '
$ws.Range('D45').Value = "'0.741"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -4.88%  '
$ws.Range('E46').Value = '  -2.23%  '
$ws.Range('D47').Value = "'53.17"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -3.79%  '
$ws.Range('D48').Value = "'0.0511"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -2.87%  '
$ws.Range('B49').Value = 'Mantle'
$ws.Range('C49').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D49').Value = "'0.410"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.03%  '
$ws.Range('B50').Value = 'RenderToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D50').Value = "'1.41"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -5.43%  '
$ws.Range('E51').Value = '  -0.03%  '
